# The deck's single slide master (ppt/theme/theme1.xml) was using the
# "Integral" design's "Red Violet" color scheme. The edit restores the
# theme to the stock Office "Office Theme" color scheme (the one that
# already lived in the unused ppt/theme/theme2.xml part, which only the
# notes master pointed to).
#
# PowerPoint's object model exposes the 12 DrawingML theme colors via
# Slide.ThemeColorScheme (MsoThemeColorSchemeIndex order: dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink). Writing .RGB on each entry rewrites
# the <a:srgbClr val="…"/> for that slot in the slide master's theme
# part in place.

function RgbVal($r, $g, $b) {
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = RgbVal 0x00 0x00 0x00   # dk1
$tcs.Item(2).RGB  = RgbVal 0xFF 0xFF 0xFF   # lt1
$tcs.Item(3).RGB  = RgbVal 0x44 0x54 0x6A   # dk2
$tcs.Item(4).RGB  = RgbVal 0xE7 0xE6 0xE6   # lt2
$tcs.Item(5).RGB  = RgbVal 0x5B 0x9B 0xD5   # accent1
$tcs.Item(6).RGB  = RgbVal 0xED 0x7D 0x31   # accent2
$tcs.Item(7).RGB  = RgbVal 0xA5 0xA5 0xA5   # accent3
$tcs.Item(8).RGB  = RgbVal 0xFF 0xC0 0x00   # accent4
$tcs.Item(9).RGB  = RgbVal 0x44 0x72 0xC4   # accent5
$tcs.Item(10).RGB = RgbVal 0x70 0xAD 0x47   # accent6
$tcs.Item(11).RGB = RgbVal 0x05 0x63 0xC1   # hlink
$tcs.Item(12).RGB = RgbVal 0x95 0x4F 0x72   # folHlink
